$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the tiny floating-point drift on the previous row's timestamp
$ws.Range("A8").Value = 45873.62529006945

# Append the new row pulled in by the scheduled task run
$ws.Range("A9").Value = 45873.70865909853
$ws.Range("B9").Value = 2025
$ws.Range("C9").Value = 15
$ws.Range("D9").Value = 19.67
$ws.Range("E9").Value = 76.20999999999999
$ws.Range("F9").Value = 145.9
$ws.Range("G9").Value = 11.9
$ws.Range("H9").Value = "ESE"
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = "17:00:28"

# Match the formatting/style applied to the other timestamp cells in column A
$ws.Range("A9").NumberFormat = $ws.Range("A8").NumberFormat
